# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5574
$ws1.Range("F5").Value = 65
$ws1.Range("F6").Value = 86
$ws1.Range("F9").Value = 527
$ws1.Range("F10").Value = 21

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5574
$ws4.Range("F6").Value = 65
$ws4.Range("F7").Value = 86
$ws4.Range("F11").Value = 527
$ws4.Range("F12").Value = 21
